$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 6

# Write the new-string-introducing cells first, in the exact order their
# values first appear so the shared-string table grows in the same order
# as the target workbook: D,E,F,G,J,K,L,N,R,Z,AE,AF,AG,AH, then A.
$ws.Cells.Item($r, 4).Value  = "MAKEPAS2713ENDOR"          # D6  MAKE
$ws.Cells.Item($r, 5).Value  = "MAKEPAS2713ENDOR"          # E6  MAKE_TEXT
$ws.Cells.Item($r, 6).Value  = "MODELPAS2713ENDOR"         # F6  MODEL_TEXT
$ws.Cells.Item($r, 7).Value  = "SERIESPAS2713ENDOR"        # G6  SERIES_TEXT
$ws.Cells.Item($r, 10).Value = "BODYTYPEPAS2713ENDOR"      # J6  BODYTYPE_TEXT
$ws.Cells.Item($r, 11).Value = "HATCHBACK 4 DOOR"          # K6  SEGMENTATION_CD
$ws.Cells.Item($r, 12).Value = "BODYSTYLEPAS2713ENDOR"     # L6  BODY_STYLE_CD
$ws.Cells.Item($r, 14).Value = "8L V12"                    # N6  ENGINE_NAME
$ws.Cells.Item($r, 18).Value = "4WD"                       # R6  WD
$ws.Cells.Item($r, 26).Value = "RT"                        # Z6  STAT
$ws.Cells.Item($r, 31).Value = "E"                         # AE6 BI_SYMBOL
$ws.Cells.Item($r, 32).Value = "E"                         # AF6 PD_SYMBOL
$ws.Cells.Item($r, 33).Value = "E"                         # AG6 UM_SYMBOL
$ws.Cells.Item($r, 34).Value = "E"                         # AH6 MP_SYMBOL
$ws.Cells.Item($r, 1).Value  = "HHHNK2CC&F"                # A6  VIN

# Remaining cells reuse existing shared strings / are plain numbers.
$ws.Cells.Item($r, 2).Value  = "SYMBOL_2000"                        # B6  VERSION
$ws.Cells.Item($r, 3).Value  = 2018                                 # C6  YEAR
$ws.Cells.Item($r, 8).Value  = 88888                                # H6  MFG_BAS_MSRP
$ws.Cells.Item($r, 9).Value  = "WAG"                                # I6  BODY
$ws.Cells.Item($r, 13).Value = "WAG"                                # M6  BODYSHELL
$ws.Cells.Item($r, 15).Value = 12                                   # O6  NUMOFCYLINDERS
$ws.Cells.Item($r, 16).Value = "G"                                  # P6  ENG_FUEL_CD
$ws.Cells.Item($r, 17).Value = 214                                  # Q6  ENG_DISPLCMNT_CI
$ws.Cells.Item($r, 19).Value = 4                                    # S6  WHEELDRIVE
$ws.Cells.Item($r, 20).Value = "000R"                               # T6  RESTRAINTSCODE
$ws.Cells.Item($r, 21).Value = "DUAL AIR BAGS FRONT"                # U6  RESTRAINTSCODE_TEXT
$ws.Cells.Item($r, 22).Value = 2                                    # V6  ANTILOCKCODE
$ws.Cells.Item($r, 23).Value = "4 WHEEL STANDARD"                   # W6  ANTILOCKCODE_TEXT
$ws.Cells.Item($r, 24).Value = "STD"                                # X6  ANTITHEFTCODE
$ws.Cells.Item($r, 25).Value = "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"   # Y6  ANTITHEFTCODE_TEXT
$ws.Cells.Item($r, 27).Value = 33                                   # AA6 COLL_SYMBOL
$ws.Cells.Item($r, 28).Value = 43                                   # AB6 COMP_SYMBOL
$ws.Cells.Item($r, 29).Value = "C"                                  # AC6 CHOICE_TIER
$ws.Cells.Item($r, 30).Value = "Y"                                  # AD6 ALTFUEL
$ws.Cells.Item($r, 35).Value = 20000101                             # AI6 ENTRYDATE
$ws.Cells.Item($r, 36).Value = "Y"                                  # AJ6 VALID
$ws.Cells.Item($r, 37).Value = "Y"                                  # AK6 ANTITHEFT_DISCOUNT
$ws.Cells.Item($r, 38).Value = "N"                                  # AL6 RESTRAINTS_DISCOUNT

# B6 reuses the same cell style already applied to B2:B5 (a plain-Calibri
# font with no theme color), so copy that formatting down instead of
# constructing a brand-new font/style entry.
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I16").Select()
